# Rename the existing sheet "Sheet 1" -> "Tab"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Tab"

# Add a new sheet called "WithTitle", placed right after "Tab"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "WithTitle"

# --- Big title row (merged A1:I1, centered) ---
$ws2.Range("A1").Value = "Big title"
$ws2.Range("A1:I1").Merge()
$ws2.Range("A1:I1").HorizontalAlignment = -4108   # xlCenter

# --- Header row (row 7): First name / Last name / Occupation, columns B:D ---
$ws2.Range("B7").Value = "First name"
$ws2.Range("C7").Value = "Last name"
$ws2.Range("D7").Value = "Occupation"

# --- Data rows 8-11 ---
$ws2.Range("B8").Value = "Michael"
$ws2.Range("C8").Value = "Jackson"
$ws2.Range("D8").Value = "Singer"

$ws2.Range("B9").Value = "Jack"
$ws2.Range("C9").Value = "The Ripper"
$ws2.Range("D9").Value = "Murderer"

$ws2.Range("B10").Value = "Stephen"
$ws2.Range("C10").Value = "King"
$ws2.Range("D10").Value = "Writer"

$ws2.Range("B11").Value = "John"
$ws2.Range("C11").Value = 3
$ws2.Range("D11").Value = "Tester"

# --- Copy cell formatting from the equivalent "Tab" cells so styles match ---
$fmtMap = @{
    "B7" = "C8";  "C7" = "D8";  "D7" = "E8";
    "B8" = "C9";  "C8" = "D9";  "D8" = "E9";
    "B9" = "C10"; "C9" = "D10"; "D9" = "E10";
    "B10" = "C11"; "C10" = "D11"; "D10" = "E11";
    "B11" = "C12"; "C11" = "D12"; "D11" = "E12";
}
foreach ($dest in $fmtMap.Keys) {
    $src = $fmtMap[$dest]
    $ws1.Range($src).Copy()
    $ws2.Range($dest).PasteSpecial(-4122)   # xlPasteFormats
}
$excel.CutCopyMode = $false

# --- Row heights for data rows 7 and 9 (25.5) ---
$ws2.Range("7:7").RowHeight = 25.5
$ws2.Range("9:9").RowHeight = 25.5

$ws2.Range("A1").Select()
